$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 97 ("「アブダビ・メディア」" post), shifting all rows below it up by one.
$ws.Rows.Item(97).Delete()
